$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D (Price) to Text format so that numeric-looking
# strings (e.g. "250.40", "0.0910") are written verbatim as text, matching
# the source data which stores these as inline/shared strings, not numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.251.59"
$ws.Range("E2").Value = "  +2.28%  "
$ws.Range("D3").Value = "2.092.31"
$ws.Range("E3").Value = "  +3.51%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "250.40"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").Value = "0.658"
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("D8").Value = "51.66"
$ws.Range("E8").Value = "  +13.64%  "
$ws.Range("D9").Value = "61.52"
$ws.Range("E9").Value = "  +4.04%  "
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("D11").Value = "0.0742"
$ws.Range("E11").Value = "  +3.35%  "
$ws.Range("E12").Value = "  +7.24%  "
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("E14").Value = "  +3.68%  "
$ws.Range("D15").Value = "0.832"
$ws.Range("E15").Value = "  +2.96%  "
$ws.Range("D16").Value = "2.098.58"
$ws.Range("E16").Value = "  +3.91%  "
$ws.Range("D17").Value = "5.11"
$ws.Range("E17").Value = "  +3.94%  "
$ws.Range("D18").Value = "37.204.40"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("D19").Value = "72.21"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("D20").Value = "14.03"
$ws.Range("E20").Value = "  +7.91%  "
$ws.Range("D21").Value = "0.0₃0838"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").Value = "240.22"
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("E23").Value = "  +6.49%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").Value = "170.54"
$ws.Range("E26").Value = "  +4.67%  "
$ws.Range("D27").Value = "9.22"
$ws.Range("E27").Value = "  +7.19%  "
$ws.Range("D28").Value = "20.69"
$ws.Range("E28").Value = "  +3.86%  "
$ws.Range("D29").Value = "2.01"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").Value = "0.122"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("B31").Value = "Gas"
$ws.Range("C31").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D31").Value = "22.85"
$ws.Range("E31").Value = "  +7.59%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "1.06"
$ws.Range("E32").Value = "  +25.81%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "4.48"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("D34").Value = "0.0609"
$ws.Range("E34").Value = "  +2.92%  "
$ws.Range("D35").Value = "0.0910"
$ws.Range("E35").Value = "  +10.59%  "
$ws.Range("D37").Value = "2.30"
$ws.Range("E37").Value = "  +7.60%  "
$ws.Range("D38").Value = "1.86"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").Value = "4.09"
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").Value = "18.30"
$ws.Range("E41").Value = "  +13.68%  "
$ws.Range("E42").Value = "  +3.77%  "
$ws.Range("E43").Value = "  +5.79%  "
$ws.Range("D44").Value = "98.54"
$ws.Range("E44").Value = "  +2.34%  "
$ws.Range("E45").Value = "  +12.70%  "
$ws.Range("D46").Value = "2.73"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.317.78"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").Value = "2.96"
$ws.Range("E48").Value = "  +6.40%  "
$ws.Range("D49").Value = "6.97"
$ws.Range("E49").Value = "  +13.28%  "
$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D50").Value = "3.77"
$ws.Range("E50").Value = "  +77.32%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.285.01"
$ws.Range("E51").Value = "  +2.43%  "

# Restore the original (default/"Normal") style on column D so no stray
# style attribute is left on the cells.
$ws.Range("D2:D51").Style = "Normal"
